# Automatic update of files.
# Rows 47 and 48 of the "Artfynd" sheet swapped their record-specific
# values (the two observations traded places) while the shared/location
# columns stayed the same.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Artfynd")

# --- Row 47 (new values = old row 48 values) ---
$ws.Range("A47").Value = 111635461
$ws.Range("B47").Value = 89590
$ws.Range("D47").Value = "VU"
$ws.Range("E47").Value = 48
$ws.Range("F47").Value = "Lappticka"
$ws.Range("G47").Value = "Amylocystis lapponica"
$ws.Range("H47").Value = "(Romell) Singer"
$ws.Range("Q47").Value = 539846.9353019162
$ws.Range("R47").Value = 7198365.604689348
$ws.Range("Z47").Value = "09:56"
$ws.Range("AB47").Value = "09:56"
$ws.Range("AW47").Value = "Yasmine Kindlund"
$ws.Range("AX47").Value = "Yasmine Kindlund, Isak Vahlström"

# --- Row 48 (new values = old row 47 values) ---
$ws.Range("A48").Value = 111634859
$ws.Range("B48").Value = 77515
$ws.Range("D48").Value = "NT"
$ws.Range("E48").Value = 6425
$ws.Range("F48").Value = "Garnlav"
$ws.Range("G48").Value = "Alectoria sarmentosa"
$ws.Range("H48").Value = "(Ach.) Ach."
$ws.Range("Q48").Value = 539847.161346367
$ws.Range("R48").Value = 7198348.622951495
$ws.Range("Z48").Value = "09:58"
$ws.Range("AB48").Value = "09:58"
$ws.Range("AW48").Value = "Isak Vahlström"
$ws.Range("AX48").Value = "Isak Vahlström, Yasmine Kindlund"
